{"js": "// Renumber the bracketed reference citations in the bibliography and\n// relocate the (Word-managed) \"_GoBack\" bookmark from the end of the\n// \"Reinforcement learning\" sentence to wrap the Watkins (1989) title.\n//\n// Mapping of old -> new leading bracket number (each occurs exactly once\n// as the start of a paragraph elsewhere in the body, so we resolve every\n// edit by locating the specific paragraph first -- this sidesteps any\n// ordering / collision issues between the renumbered values):\n//   [4]  -> [7]\n//   [5]  -> []     (number dropped)\n//   [6]  -> [4]\n//   [7]  -> [5]\n//   [8]  -> [6]\n//   [9]  -> [8]\n//   [10] -> [9]\n//   [11] -> [10]\n//   [12] -> [12]   (unchanged)\n//   [13] -> [11]\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Replace the FIRST occurrence of `find` inside `paragraph` with `replace`,\n// matching case-sensitively and leaving surrounding runs/formatting intact.\nfunction renumber(paragraph, find, replace) {\n  const results = paragraph.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  context.trackedObjects.add(results);\n  return results;\n}\n\nconst searches = [];\nconst items = paragraphs.items;\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (t.indexOf(\"[4] \") === 0) searches.push([renumber(items[i], \"[4]\"), \"[7]\"]);\n  else if (t.indexOf(\"[5]\") === 0) searches.push([renumber(items[i], \"[5]\"), \"[]\"]);\n  else if (t.indexOf(\"[6] \") === 0) searches.push([renumber(items[i], \"[6]\"), \"[4]\"]);\n  else if (t.indexOf(\"[7]\") === 0) searches.push([renumber(items[i], \"[7]\"), \"[5]\"]);\n  else if (t.indexOf(\"[8] Watkins\") === 0) searches.push([renumber(items[i], \"[8]\"), \"[6]\"]);\n  else if (t.indexOf(\"[9] \") === 0) searches.push([renumber(items[i], \"[9]\"), \"[8]\"]);\n  else if (t.indexOf(\"[10] \") === 0) searches.push([renumber(items[i], \"[10]\"), \"[9]\"]);\n  else if (t.indexOf(\"[11]\") === 0) searches.push([renumber(items[i], \"[11]\"), \"[10]\"]);\n  else if (t.indexOf(\"[13] \") === 0) searches.push([renumber(items[i], \"[13]\"), \"[11]\"]);\n}\n\nawait context.sync();\n\nfor (const [results, newText] of searches) {\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// Relocate the \"_GoBack\" bookmark: delete it from its current position\n// (right after the comma that follows \"Reinforcement learning\") and\n// re-insert it around the Watkins reference title (excluding the final\n// period), matching the new location in the edited document.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst watkinsResults = body.search(\n  \"Watkins, C. J. C. H. (1989). Learning from delayed rewards\",\n  { matchCase: true }\n);\nwatkinsResults.load(\"text\");\nawait context.sync();\n\nwatkinsResults.items[0].insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Renumber the bracketed reference citations in the bibliography and\n# relocate the (Word-managed) \"_GoBack\" bookmark from the end of the\n# \"Reinforcement learning\" sentence to wrap the Watkins (1989) title.\n#\n# Mapping of old -> new leading bracket number (each bracket number occurs\n# exactly once, as the start of a specific paragraph, elsewhere in the\n# body -- so every substitution is scoped to its own paragraph's Range,\n# which makes the order irrelevant and avoids any collisions between the\n# renumbered values):\n#   [4]  -> [7]\n#   [5]  -> []     (number dropped)\n#   [6]  -> [4]\n#   [7]  -> [5]\n#   [8]  -> [6]\n#   [9]  -> [8]\n#   [10] -> [9]\n#   [11] -> [10]\n#   [13] -> [11]\n\n$d = $word.ActiveDocument\n\n# NOTE: this interpreter only reliably binds POSITIONAL parameters, so the\n# helper below is always called as `Renumber-Paragraph <idx> <find> <repl>`.\nfunction Renumber-Paragraph {\n    param([int]$ParaIndex, [string]$FindText, [string]$ReplaceText)\n    $r = $d.Paragraphs($ParaIndex).Range\n    $r.Find.ClearFormatting()\n    $r.Find.Execute($FindText, $false, $false, $false, $false, $false, $true, 1, $false, $ReplaceText, 2) | Out-Null\n}\n\nRenumber-Paragraph 4  \"[4]\"  \"[7]\"\nRenumber-Paragraph 6  \"[5]\"  \"[]\"\nRenumber-Paragraph 7  \"[6]\"  \"[4]\"\nRenumber-Paragraph 8  \"[7]\"  \"[5]\"\nRenumber-Paragraph 9  \"[8]\"  \"[6]\"\nRenumber-Paragraph 10 \"[9]\"  \"[8]\"\nRenumber-Paragraph 12 \"[10]\" \"[9]\"\nRenumber-Paragraph 14 \"[11]\" \"[10]\"\nRenumber-Paragraph 18 \"[13]\" \"[11]\"\n\n# Relocate the \"_GoBack\" bookmark: find the Watkins reference sentence\n# (excluding the trailing period) and re-add the bookmark there; adding a\n# bookmark with an already-existing name moves it instead of creating a\n# duplicate.\n$watkinsPara = $d.Paragraphs(9).Range\n$watkinsPara.Find.ClearFormatting()\n$watkinsPara.Find.Execute(\"Watkins, C. J. C. H. (1989). Learning from delayed rewards\") | Out-Null\n$d.Bookmarks.Add(\"_GoBack\", $watkinsPara) | Out-Null\n"}
